$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Rename the first sheet
$ws1.Name = "Lista ansökningar"

# Shrink the auto-fit columns A-C to narrower, explicit widths
# (values chosen so the stored OOXML column width lands as close as
# possible to the target widths of 16.7109375 / 60.5703125 / 23.5703125)
$ws1.Columns.Item(1).ColumnWidth = 15.833333333333334
$ws1.Columns.Item(2).ColumnWidth = 59.666666666666664
$ws1.Columns.Item(3).ColumnWidth = 22.666666666666668

# Move the active selection on sheet 1 to G1
$ws1.Activate() | Out-Null
$ws1.Range("G1").Select() | Out-Null
